$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added at row 131, pushing the existing
# rows 131-204 down to 132-205 (the former row 204 ends up duplicated
# as the new row 205).
$ws.Rows.Item(131).Insert()

$ws.Range("A131").Value = 2
$ws.Range("B131").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C131").Value = "Coquimbo"
$ws.Range("D131").Value = 45001
$ws.Range("E131").Value = 4
$ws.Range("F131").Value = 100112043
$ws.Range("G131").Value = "Pepino ensalada"
$ws.Range("H131").Value = "Sin especificar"
$ws.Range("I131").Value = "Primera"
$ws.Range("J131").Value = 500
$ws.Range("K131").Value = 7000
$ws.Range("L131").Value = 8000
$ws.Range("M131").Value = 7500
$ws.Range("N131").Value = "`$/caja 70 unidades"
$ws.Range("O131").Value = "Provincia de Limarí"
$ws.Range("P131").Value = 107
$ws.Range("Q131").Value = 70
$ws.Range("R131").Value = "Hortaliza"
